$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Sheet1")
$ws2 = $wb.Worksheets.Item("Sheet2")

# ------------------------------------------------------------------
# Sheet2: remove the feeds that are now configured on Sheet1.
# Rows 1:7 (SpongeIron, PigIron, Ingots, SteelOpenHearthFurnaces,
# LongRolledProducts, FlatRolledProducts) already exist on Sheet1, so
# they disappear from the "still to configure" list on Sheet2.
# ------------------------------------------------------------------
$ws2.Rows("1:7").ClearContents()

# Ids 9-13 (HotRolledProducts, ContinuouslyCastSteel,
# LiquidSteelForCastings, TotalProductionOfCrudeSteel,
# RailwayTrackMaterial) lose their id/name - four of them move to
# Sheet1 below, and RailwayTrackMaterial is dropped altogether.
$ws2.Range("A8:B12").ClearContents()

# ------------------------------------------------------------------
# Sheet1: add the four newly configured feeds.
# ------------------------------------------------------------------
$ws1.Range("A9:E9").Copy()
$ws1.Range("A10:E12").PasteSpecial(-4122)   # xlPasteFormats

$ws1.Range("A10").Value = 12
$ws1.Range("B10").Value = "TotalProductionOfCrudeSteel"
$ws1.Range("C10").Value = "WSACrudeSteel"
$ws1.Range("D10").Value = "25-32"
$ws1.Range("E10").Value = "35-44"

$ws1.Range("A11").Value = 10
$ws1.Range("B11").Value = "ContinuouslyCastSteel"
$ws1.Range("C11").Value = "WSAContinuouslyCastSteel"
$ws1.Range("D11").Value = "28-32"
$ws1.Range("E11").Value = "39-44"

$ws1.Range("A12").Value = 11
$ws1.Range("B12").Value = "LiquidSteelForCastings"
$ws1.Range("C12").Value = "WSALiquidSteelForCastings"
$ws1.Range("D12").Value = "20-24"
$ws1.Range("E12").Value = "30-34"

$ws1.Range("A9:E9").Copy()
$ws1.Range("A13:E13").PasteSpecial(-4122)   # xlPasteFormats
$ws1.Range("C13").Style = "Normal"

$ws1.Range("A13").Value = 9
$ws1.Range("B13").Value = "HotRolledProducts"
$ws1.Range("C13").Value = "WSAHotRolledProducts"
$ws1.Range("D13").Value = "33-63"
$ws1.Range("E13").Value = "45-231"

# Trailing (empty) row left below the table, matching the unstyled
# cells produced when the blank formatting of an untouched cell is
# pasted over B14:C14.
$ws1.Range("F1").Copy()
$ws1.Range("B14:C14").PasteSpecial(-4122)   # xlPasteFormats
$excel.CutCopyMode = $false

# Column C grew wider to fit the new, longer feed names.
$ws1.Columns.Item(3).AutoFit()

# ------------------------------------------------------------------
# View/selection state
# ------------------------------------------------------------------
$ws2.Range("A12:B12").Select()
$ws2.Columns.Item(2).AutoFit()
try { $excel.ActiveWindow.ScrollRow = 6 } catch {}

$ws1.Range("B1").Select()
